# Development Update | ConPass Arrival & Assignment Messages
#
# This script reproduces, via Excel COM-interop, the edits described by the
# commit: several "ConPass" / "Named Pipe" TODO markers are resolved to
# "Yes", two stale remark strings are removed (one of them replaced by a new
# remark), a new note is added, and the active sheet/selection changes from
# ConPass to CenTaxi.

$wb = $excel.ActiveWorkbook

$wsCenDLL  = $wb.Worksheets.Item("CenDLL")
$wsCTDLL   = $wb.Worksheets.Item("CTDLL")
$wsCenTaxi = $wb.Worksheets.Item("CenTaxi")
$wsConTaxi = $wb.Worksheets.Item("ConTaxi")
$wsConPass = $wb.Worksheets.Item("ConPass")

# ----------------------------------------------------------------------
# CTDLL: remark "Needs completion" on the "assigned passenger" row is
# resolved/removed.
# ----------------------------------------------------------------------
$wsCTDLL.Range("D9").Value = ""

# ----------------------------------------------------------------------
# CenTaxi: a handful of "ConPass" remarks are resolved to "Yes", the
# "Needs to be looked at" remark is cleared, and one "ConPass" remark is
# replaced with a new note about movement being implemented first.
# ----------------------------------------------------------------------
$wsCenTaxi.Range("D6").Value = ""

$wsCenTaxi.Range("C22").Value = "Yes"
$wsCenTaxi.Range("D22").Value = ""

$wsCenTaxi.Range("D23").Value = "After movement is implemented"

$wsCenTaxi.Range("C27").Value = "Yes"
$wsCenTaxi.Range("D27").Value = ""

$wsCenTaxi.Range("C30").Value = "Yes"
$wsCenTaxi.Range("D30").Value = ""

# ----------------------------------------------------------------------
# ConTaxi: "Named Pipe" remark is resolved to "Yes".
# ----------------------------------------------------------------------
$wsConTaxi.Range("C11").Value = "Yes"
$wsConTaxi.Range("D11").Value = ""

# ----------------------------------------------------------------------
# ConPass: resolve a couple of items to "Yes", add a new D8 note, drop the
# stale "Missing recover" note, and reword the assumption about passenger
# removal.
# ----------------------------------------------------------------------
$wsConPass.Range("C7").Value = "Yes"

$wsConPass.Range("D8").Value = "After movement is implemented"

$wsConPass.Range("D9").Value = ""

$wsConPass.Range("E11").Value = "- Passengers are not removed by own choice"

$wsConPass.Range("C12").Value = "Yes"

# ----------------------------------------------------------------------
# View-state: the active sheet moves from ConPass to CenTaxi, and the
# remembered selection on several sheets changes. Select on each sheet in
# turn (this also activates it) so every sheet's stored selection updates,
# finishing on CenTaxi so it ends up the active tab.
# ----------------------------------------------------------------------
$wsCTDLL.Range("D14").Select()
$wsConPass.Range("E16").Select()
$wsCenTaxi.Range("E29").Select()
